# Update DM integration fixture hierarchies
#
# Replaces the UUID identifiers (column A) in the four worksheets of the
# "explicit-dimensions-2018-1" fixture workbook, and re-applies the
# column-A autofit width that results from the new identifier values on
# the "Codes" (sheet2), "Extensions" (sheet3) and "Members_dpmDimension"
# (sheet4) sheets. ("CodeSchemes", sheet1, keeps its original column-A
# width - it is not touched by this change.)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: CodeSchemes
# ---------------------------------------------------------------------
$wsCodeSchemes = $wb.Worksheets.Item("CodeSchemes")
$wsCodeSchemes.Range("A2").Value = "ea439bf8-4eec-4df5-916a-90b144e1d58b"

# ---------------------------------------------------------------------
# Sheet: Codes
# ---------------------------------------------------------------------
$wsCodes = $wb.Worksheets.Item("Codes")
$wsCodes.Range("A2").Value = "b042fb98-44f3-41b4-b033-8ad4e766bcf4"
$wsCodes.Range("A3").Value = "2bcedb0f-6fe3-41d8-a887-c53337b9ce71"
$wsCodes.Range("A4").Value = "63cd4d2e-cdde-4abb-ae29-6f575b4addbc"
$wsCodes.Range("A5").Value = "c6334b1d-c45e-42e8-9bd9-984576e388bc"
$wsCodes.Columns.Item(1).ColumnWidth = 34.42857142857143

# ---------------------------------------------------------------------
# Sheet: Extensions
# ---------------------------------------------------------------------
$wsExtensions = $wb.Worksheets.Item("Extensions")
$wsExtensions.Range("A2").Value = "61cd9133-85bf-49df-a4cc-b02f4b7e0cd4"
$wsExtensions.Columns.Item(1).ColumnWidth = 29.0

# ---------------------------------------------------------------------
# Sheet: Members_dpmDimension
# ---------------------------------------------------------------------
$wsMembers = $wb.Worksheets.Item("Members_dpmDimension")
$wsMembers.Range("A2").Value = "e363b4d7-504f-4a4b-bed6-41bdb1f03659"
$wsMembers.Range("A3").Value = "b03d2cd4-4416-429a-bdaa-e661f7241d30"
$wsMembers.Range("A4").Value = "7e27e7ca-9ca2-401a-9955-40eea1d8c443"
$wsMembers.Range("A5").Value = "735c9b13-795a-4328-bba3-7d2b15832476"
$wsMembers.Columns.Item(1).ColumnWidth = 34.42857142857143
